# Update the weekly fruit/vegetable price dataset (Caqui - Agro Chillan)
# Rows 2-7 are refreshed with a new week's rotation of values for
# columns D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de
# comercializacion), R (Origen), S (Precio $/Kg) and T (Kg / unidad).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D = 44698; L = "Primera"; M = 120; N = 16000; O = 17000; P = 16500; Q = "`$/caja 18 kilos granel"; R = "Región de O'Higgins"; S = 917;  T = 18 }
    3 = @{ D = 44334; L = "Primera"; M = 120; N = 12000; O = 13000; P = 12500; Q = "`$/caja 12 kilos empedrada"; R = "Región de O'Higgins"; S = 1042; T = 12 }
    4 = @{ D = 44344; L = "Primera"; M = 120; N = 13000; O = 14000; P = 13500; Q = "`$/caja 18 kilos granel"; R = "Provincia de Curicó"; S = 750;  T = 18 }
    5 = @{ D = 44330; L = "Primera"; M = 60;  N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos granel"; R = "Provincia de Curicó"; S = 861;  T = 18 }
    6 = @{ D = 44316; L = "Primera"; M = 60;  N = 17500; O = 18000; P = 17750; Q = "`$/caja 16 kilos granel"; R = "Región de O'Higgins"; S = 1109; T = 16 }
    7 = @{ D = 44316; L = "Segunda"; M = 40;  N = 16000; O = 16000; P = 16000; Q = "`$/caja 16 kilos granel"; R = "Región de O'Higgins"; S = 1000; T = 16 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
    $ws.Range("S$row").Value = $vals.S
    $ws.Range("T$row").Value = $vals.T
}
